# Update cryptocurrency Price (D) and Volume(1h) (E) columns on Sheet1
# to reflect the latest coinranking.com snapshot (GitHub Actions refresh).
#
# Note: some Price values are plain decimals (e.g. "216.84") that Excel would
# otherwise auto-convert to a Number on assignment. The source data keeps the
# Price column as text (it also holds thousand-dotted values like "27.110.91"
# that cannot be numbers), so we prefix those plain-decimal values with a
# leading apostrophe - exactly what typing them into Excel by hand would do -
# to force them to stay text instead of becoming numeric.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.110.91"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "1.638.41"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'216.84"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Value = "'0.519"
$ws.Range("E6").Value = "  +1.96%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").Value = "1.867.39"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").Value = "1.649.33"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "'4.12"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("E15").Value = "  +2.36%  "

$ws.Range("E16").Value = "  -0.56%  "

$ws.Range("D17").Value = "27.110.72"
$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("E18").Value = "  +1.39%  "

$ws.Range("D19").Value = "'217.33"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("E21").Value = "  +1.88%  "

$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("E23").Value = "  +3.66%  "

$ws.Range("D24").Value = "'9.08"
$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("D25").Value = "'146.29"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +1.49%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("E30").Value = "  +1.21%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  +1.67%  "

$ws.Range("E33").Value = "  +0.84%  "

$ws.Range("D34").Value = "1.308.13"
$ws.Range("E34").Value = "  +3.34%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("E36").Value = "  +1.17%  "

$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("E38").Value = "  +2.99%  "

$ws.Range("D39").Value = "'0.543"
$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("D41").Value = "'0.809"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("E42").Value = "  +5.47%  "

$ws.Range("D43").Value = "'5.30"
$ws.Range("E43").Value = "  -1.47%  "

$ws.Range("D44").Value = "1.777.48"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "'91.38"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D50").Value = "'7.61"
$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").Value = "'0.0959"
$ws.Range("E51").Value = "  +0.11%  "
